# Bump the "Förändrad" (Changed) date column (C) from 2025-05-07 (45784)
# to 2025-05-08 (45785) for every data row (rows 2 through 43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C43").Value = 45785
